$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the current page bookmark for "Researching Information Systems and
# Computing" (row 11) from 141 to 144 pages read.
$ws.Range("C11").Value = 144

# Move the active selection to C14, matching where the user clicked next.
$ws.Range("C14").Select()

$wb.Save()
